$wb = $excel.ActiveWorkbook

# "Overview" sheet: G2 = Latest HO Xliff Generate Date for the 65ba8943... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 13:04:43"

# "zh-cn" sheet: H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 13:04:38"
$wsZhCn.Range("K2").Value = "2016-08-16 13:04:57"

# "de-de" sheet: K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-16 13:05:17"
